# Simulated Wild Card round and logged it
# Update Home row (row 2) target depth data on both the OFF and DEF sheets
# to reflect the additional game's stats.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 346
$wsOff.Range("C2").Value = 243
$wsOff.Range("D2").Value = 79
$wsOff.Range("E2").Value = 33
$wsOff.Range("F2").Value = 6

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 339
$wsDef.Range("C2").Value = 239
$wsDef.Range("D2").Value = 71
$wsDef.Range("E2").Value = 22
$wsDef.Range("F2").Value = 8
